$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits right
#    under the H1 title at the top of the document.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Meta description:*") {
        $para.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Insert a new bold paragraph ("Play Double Lucky Mushrooms DoubleMax for
#    Free") right before the paragraph that holds the "Create a feature
#    image..." image prompt (the last paragraph in the document).
#    We build it via InsertXML so the run layout (a leading empty <w:r/>
#    followed by the bold run) matches the rest of the document's style.
#    A throw-away marker paragraph is appended in the same call so that the
#    insertion point's own paragraph break is produced correctly; the
#    marker text is deleted right after, leaving the original trailing
#    paragraph intact immediately after our new paragraph.
# ---------------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Create a feature image*") {
        $targetPara = $para
        break
    }
}

if ($targetPara -ne $null) {
    $precedingEnd = $targetPara.Range.Start
    $insertionPoint = $d.Range($precedingEnd, $precedingEnd)
    $insertStart = $insertionPoint.Start

    $boldText = "Play Double Lucky Mushrooms DoubleMax for Free"
    $marker = "MARKERZZZ"
    $xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>' + $boldText + '</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>' + $marker + '</w:t></w:r></w:p>'
    $insertionPoint.InsertXML($xmlFrag)

    $markerStart = $insertStart + $boldText.Length + 1
    $markerEnd = $markerStart + $marker.Length
    $markerRange = $d.Range($markerStart, $markerEnd)
    $markerRange.Delete()
}

# ---------------------------------------------------------------------------
# 3) Replace the old image-prompt text with the meta-description copy,
#    keeping its italic run formatting untouched.
# ---------------------------------------------------------------------------
$oldText = "Create a feature image that perfectly captures the excitement and magic of Double Lucky Mushrooms DoubleMax. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be surrounded by the lush green fields of Ireland, with a rainbow in clear sight behind them. Include elements from the game, such as the expanding Wild clover, pot of gold coins, and mushroom symbols, to give players a glimpse of what they can expect from the game. Use bright and bold colors to make the image stand out and grab players' attention."
$newText = "Read a review of Double Lucky Mushrooms DoubleMax, a high volatility slot game with beautiful graphics and try it out for free on any device."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
